$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 with same style as other header cells (copy format from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Data values for F2:F86
$timeTaken = @(
  "2021-10-05 13:39:03.349854",
  "2021-10-05 13:39:03.349865",
  "2021-10-05 13:39:03.349869",
  "2021-10-05 13:39:03.349872",
  "2021-10-05 13:39:03.349875",
  "2021-10-05 13:39:03.349878",
  "2021-10-05 13:39:03.349880",
  "2021-10-05 13:39:03.349883",
  "2021-10-05 13:39:03.349886",
  "2021-10-05 13:39:03.349889",
  "2021-10-05 13:39:03.349891",
  "2021-10-05 13:39:03.349894",
  "2021-10-05 13:39:03.349897",
  "2021-10-05 13:39:03.349900",
  "2021-10-05 13:39:03.349902",
  "2021-10-05 13:39:03.349905",
  "2021-10-05 13:39:03.349908",
  "2021-10-05 13:39:03.349911",
  "2021-10-05 13:39:03.349913",
  "2021-10-05 13:39:03.349916",
  "2021-10-05 13:39:03.349919",
  "2021-10-05 13:39:03.349922",
  "2021-10-05 13:39:03.349924",
  "2021-10-05 13:39:03.349927",
  "2021-10-05 13:39:03.349930",
  "2021-10-05 13:39:03.349933",
  "2021-10-05 13:39:03.349935",
  "2021-10-05 13:39:03.349938",
  "2021-10-05 13:39:03.349941",
  "2021-10-05 13:39:03.349943",
  "2021-10-05 13:39:03.349947",
  "2021-10-05 13:39:03.349949",
  "2021-10-05 13:39:03.349952",
  "2021-10-05 13:39:03.349955",
  "2021-10-05 13:39:03.349958",
  "2021-10-05 13:39:03.349961",
  "2021-10-05 13:39:03.349964",
  "2021-10-05 13:39:03.349966",
  "2021-10-05 13:39:03.349969",
  "2021-10-05 13:39:03.349972",
  "2021-10-05 13:39:03.349975",
  "2021-10-05 13:39:03.349978",
  "2021-10-05 13:39:03.349981",
  "2021-10-05 13:39:03.349984",
  "2021-10-05 13:39:03.349987",
  "2021-10-05 13:39:03.349989",
  "2021-10-05 13:39:03.349992",
  "2021-10-05 13:39:03.349995",
  "2021-10-05 13:39:03.349998",
  "2021-10-05 13:39:03.350000",
  "2021-10-05 13:39:03.350003",
  "2021-10-05 13:39:03.350006",
  "2021-10-05 13:39:03.350009",
  "2021-10-05 13:39:03.350012",
  "2021-10-05 13:39:03.350015",
  "2021-10-05 13:39:03.350018",
  "2021-10-05 13:39:03.350020",
  "2021-10-05 13:39:03.350023",
  "2021-10-05 13:39:03.350026",
  "2021-10-05 13:39:03.350028",
  "2021-10-05 13:39:03.350031",
  "2021-10-05 13:39:03.350034",
  "2021-10-05 13:39:03.350036",
  "2021-10-05 13:39:03.350039",
  "2021-10-05 13:39:03.350043",
  "2021-10-05 13:39:03.350047",
  "2021-10-05 13:39:03.350050",
  "2021-10-05 13:39:03.350052",
  "2021-10-05 13:39:03.350055",
  "2021-10-05 13:39:03.350058",
  "2021-10-05 13:39:03.350061",
  "2021-10-05 13:39:03.350063",
  "2021-10-05 13:39:03.350066",
  "2021-10-05 13:39:03.350069",
  "2021-10-05 13:39:03.350071",
  "2021-10-05 13:39:03.350074",
  "2021-10-05 13:39:03.350079",
  "2021-10-05 13:39:03.350083",
  "2021-10-05 13:39:03.350085",
  "2021-10-05 13:39:03.350088",
  "2021-10-05 13:39:03.350091",
  "2021-10-05 13:39:03.350094",
  "2021-10-05 13:39:03.350096",
  "2021-10-05 13:39:03.350099",
  "2021-10-05 13:39:03.350102"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
  $ws.Cells.Item($i + 2, 6).Value = $timeTaken[$i]
}

